$wb = $excel.ActiveWorkbook

# Update 2024 (column K, and a couple 2023/column J corrections) figures
# for "2024-05-12" data refresh across the Citywide Totals, By Neighborhood
# summary, and individual neighborhood sheets.

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2633
$ws.Range("K3").Value = 2546
$ws.Range("K4").Value = 534
$ws.Range("K5").Value = 169
$ws.Range("K6").Value = 3173
$ws.Range("K7").Value = 9055

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("K5").Value = 4
$ws.Range("K6").Value = 9

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 27
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 175
$ws.Range("K6").Value = 200
$ws.Range("K7").Value = 600

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 76
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 99
$ws.Range("K3").Value = 129
$ws.Range("K6").Value = 100
$ws.Range("K7").Value = 356

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 80
$ws.Range("K4").Value = 13
$ws.Range("K5").Value = 14
$ws.Range("K7").Value = 299

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 62
$ws.Range("K7").Value = 212

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K8").Value = 600
$ws.Range("K10").Value = 51
$ws.Range("K11").Value = 191
$ws.Range("K14").Value = 55
$ws.Range("K17").Value = 16
$ws.Range("K19").Value = 265
$ws.Range("K20").Value = 207
$ws.Range("K23").Value = 80
$ws.Range("K25").Value = 36
$ws.Range("K27").Value = 98
$ws.Range("K29").Value = 467
$ws.Range("K32").Value = 15
$ws.Range("K33").Value = 356
$ws.Range("K34").Value = 45
$ws.Range("K37").Value = 299
$ws.Range("K38").Value = 9
$ws.Range("K41").Value = 76
$ws.Range("J42").Value = 1239
$ws.Range("K42").Value = 312
$ws.Range("K47").Value = 50
$ws.Range("K51").Value = 99
$ws.Range("K52").Value = 251
$ws.Range("K53").Value = 133
$ws.Range("K54").Value = 167
$ws.Range("K55").Value = 99
$ws.Range("K56").Value = 12
$ws.Range("J63").Value = 98
$ws.Range("K63").Value = 32
$ws.Range("K65").Value = 212
$ws.Range("K66").Value = 32
$ws.Range("K67").Value = 351
$ws.Range("K71").Value = 27
$ws.Range("K73").Value = 90
$ws.Range("K76").Value = 135
$ws.Range("K78").Value = 126
$ws.Range("K83").Value = 199
$ws.Range("K85").Value = 432
$ws.Range("K89").Value = 119
$ws.Range("K94").Value = 105
$ws.Range("K96").Value = 129
$ws.Range("K101").Value = 9055

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 351

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 167

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 157
$ws.Range("K4").Value = 28
$ws.Range("K6").Value = 148
$ws.Range("K7").Value = 467

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 71
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 23
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 79
$ws.Range("K3").Value = 99
$ws.Range("J6").Value = 654
$ws.Range("J7").Value = 1239
$ws.Range("K7").Value = 312

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 23
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 24
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 72
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 207

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K5").Value = 8
$ws.Range("K7").Value = 432

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 71
$ws.Range("K4").Value = 13
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 251
